# Auto-generated edit script applying the Raiden_Profits.xlsx diff
# Updates market-price-derived columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 344.75
$ws.Range("I28").Value = 363
$ws.Range("J28").Value = 290
$ws.Range("K28").Value = 363
$ws.Range("L28").Value = 290
$ws.Range("M28").Value = 122
$ws.Range("N28").Value = -1260
$ws.Range("H33").Value = 334.08694
$ws.Range("I33").Value = 232.4762
$ws.Range("K33").Value = 232.4762
$ws.Range("M33").Value = -3.476200000000006
$ws.Range("H63").Value = 69000
$ws.Range("J63").Value = 69000
$ws.Range("L63").Value = 69000
$ws.Range("N63").Value = -70248
$ws.Range("H66").Value = 69000
$ws.Range("J66").Value = 69000
$ws.Range("L66").Value = 207000
$ws.Range("N66").Value = -213240
$ws.Range("H70").Value = 156255.16
$ws.Range("I70").Value = 242129.17
$ws.Range("J70").Value = 9042.571
$ws.Range("K70").Value = 726387.51
$ws.Range("L70").Value = 27127.713
$ws.Range("M70").Value = -726117.51
$ws.Range("N70").Value = -27667.713
$ws.Range("H73").Value = 156255.16
$ws.Range("I73").Value = 242129.17
$ws.Range("J73").Value = 9042.571
$ws.Range("K73").Value = 726387.51
$ws.Range("L73").Value = 27127.713
$ws.Range("M73").Value = -725451.51
$ws.Range("N73").Value = -28999.713
$ws.Range("H74").Value = 4342
$ws.Range("I74").Value = 4342
$ws.Range("K74").Value = 4342
$ws.Range("M74").Value = -3406
$ws.Range("H75").Value = 73439.28999999999
$ws.Range("J75").Value = 73439.28999999999
$ws.Range("L75").Value = 73439.28999999999
$ws.Range("N75").Value = -75311.28999999999
$ws.Range("H76").Value = 8689.458000000001
$ws.Range("I76").Value = 8327.5
$ws.Range("K76").Value = 8327.5
$ws.Range("M76").Value = -8012.5
$ws.Range("H77").Value = 4342
$ws.Range("I77").Value = 4342
$ws.Range("K77").Value = 21710
$ws.Range("M77").Value = -17030
$ws.Range("H78").Value = 73439.28999999999
$ws.Range("J78").Value = 73439.28999999999
$ws.Range("L78").Value = 220317.87
$ws.Range("N78").Value = -229677.87
$ws.Range("H79").Value = 8689.458000000001
$ws.Range("I79").Value = 8327.5
$ws.Range("K79").Value = 8327.5
$ws.Range("M79").Value = -7235.5
$ws.Range("H80").Value = 559.25
$ws.Range("J80").Value = 502.77777
$ws.Range("L80").Value = 1508.33331
$ws.Range("N80").Value = -3504.33331
$ws.Range("H81").Value = 45000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996
$ws.Range("H83").Value = 559.25
$ws.Range("J83").Value = 502.77777
$ws.Range("L83").Value = 4524.99993
$ws.Range("N83").Value = -14508.99993
$ws.Range("H84").Value = 45000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984
$ws.Range("H99").Value = 1078.4
$ws.Range("I99").Value = 1097.5
$ws.Range("J99").Value = 1065.6666
$ws.Range("K99").Value = 3292.5
$ws.Range("L99").Value = 3196.9998
$ws.Range("M99").Value = -1794.5
$ws.Range("N99").Value = -6192.9998
$ws.Range("H100").Value = 2247
$ws.Range("I100").Value = 2033.3334
$ws.Range("K100").Value = 2033.3334
$ws.Range("M100").Value = -1492.3334
$ws.Range("H111").Value = 2971.3333
$ws.Range("I111").Value = 2925.2222
$ws.Range("J111").Value = 3017.4443
$ws.Range("K111").Value = 8775.6666
$ws.Range("L111").Value = 9052.332900000001
$ws.Range("M111").Value = -5708.6666
$ws.Range("N111").Value = -15186.3329
$ws.Range("H116").Value = 3491.8125
$ws.Range("I116").Value = 3339.75
$ws.Range("K116").Value = 3339.75
$ws.Range("M116").Value = 102.25
$ws.Range("H125").Value = 3250
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 13500
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -11040
$ws.Range("N125").Value = -49920
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 67500
$ws.Range("J130").Value = 67500
$ws.Range("L130").Value = 67500
$ws.Range("N130").Value = -77540
$ws.Range("H131").Value = 1706.8334
$ws.Range("I131").Value = 1061
$ws.Range("K131").Value = 3183
$ws.Range("M131").Value = 1857
$ws.Range("H132").Value = 385877.7
$ws.Range("I132").Value = 1354.75
$ws.Range("K132").Value = 4064.25
$ws.Range("M132").Value = -1534.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 1259.7142
$ws.Range("I135").Value = 1445.4
$ws.Range("J135").Value = 795.5
$ws.Range("K135").Value = 13008.6
$ws.Range("L135").Value = 7159.5
$ws.Range("M135").Value = -10473.6
$ws.Range("N135").Value = -12229.5
$ws.Range("H136").Value = 69999
$ws.Range("J136").Value = 69999
$ws.Range("L136").Value = 69999
$ws.Range("N136").Value = -80199
$ws.Range("H137").Value = 4294.3335
$ws.Range("I137").Value = 2070.3
$ws.Range("K137").Value = 6210.900000000001
$ws.Range("M137").Value = -3660.900000000001
$ws.Range("H139").Value = 67500
$ws.Range("J139").Value = 67500
$ws.Range("L139").Value = 67500
$ws.Range("N139").Value = -77780
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14558.589
$ws.Range("I32").Value = 2569.628
$ws.Range("K32").Value = 2569.628
$ws.Range("M32").Value = -2282.628
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4480
$ws.Range("N39").ClearContents()
$ws.Range("H88").Value = 1643
$ws.Range("I88").Value = 1004.8333
$ws.Range("K88").Value = 1004.8333
$ws.Range("M88").Value = -598.8333
$ws.Range("H91").Value = 1643
$ws.Range("I91").Value = 1004.8333
$ws.Range("K91").Value = 1004.8333
$ws.Range("M91").Value = 399.1667
$ws.Range("H94").Value = 29998.5
$ws.Range("J94").Value = 29998.5
$ws.Range("L94").Value = 29998.5
$ws.Range("N94").Value = -31800.5
$ws.Range("H97").Value = 1009.2692
$ws.Range("I97").Value = 1034.2
$ws.Range("K97").Value = 1034.2
$ws.Range("M97").Value = -538.2
$ws.Range("H110").Value = 2957
$ws.Range("I110").Value = 2460.5
$ws.Range("J110").Value = 4446.5
$ws.Range("K110").Value = 2460.5
$ws.Range("L110").Value = 4446.5
$ws.Range("M110").Value = -415.5
$ws.Range("N110").Value = -8536.5
$ws.Range("H122").Value = 2224.6428
$ws.Range("I122").Value = 2082.9697
$ws.Range("J122").Value = 2744.111
$ws.Range("K122").Value = 6248.909100000001
$ws.Range("L122").Value = 8232.332999999999
$ws.Range("M122").Value = -3798.909100000001
$ws.Range("N122").Value = -13132.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1141.4814
$ws.Range("I107").Value = 1118.4584
$ws.Range("J107").Value = 1325.6666
$ws.Range("K107").Value = 1118.4584
$ws.Range("L107").Value = 1325.6666
$ws.Range("M107").Value = 801.5416
$ws.Range("N107").Value = -5165.6666
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5157.14
$ws.Range("I31").Value = 3330.1052
$ws.Range("K31").Value = 3330.1052
$ws.Range("M31").Value = -3035.1052
$ws.Range("H34").Value = 5157.14
$ws.Range("I34").Value = 3330.1052
$ws.Range("K34").Value = 3330.1052
$ws.Range("M34").Value = -3128.1052
$ws.Range("H52").Value = 90633
$ws.Range("J52").Value = 90633
$ws.Range("L52").Value = 90633
$ws.Range("N52").Value = -91221
$ws.Range("H58").Value = 3350.0454
$ws.Range("J58").Value = 1672
$ws.Range("L58").Value = 1672
$ws.Range("N58").Value = -2078
$ws.Range("H62").Value = 12321.667
$ws.Range("I62").Value = 12321.667
$ws.Range("K62").Value = 12321.667
$ws.Range("M62").Value = -11697.667
$ws.Range("H65").Value = 12321.667
$ws.Range("I65").Value = 12321.667
$ws.Range("K65").Value = 61608.335
$ws.Range("M65").Value = -58488.335
$ws.Range("H106").Value = 39999.5
$ws.Range("J106").Value = 39999.5
$ws.Range("L106").Value = 39999.5
$ws.Range("N106").Value = -42523.5
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H132").Value = 3186.5454
$ws.Range("I132").Value = 1881.625
$ws.Range("K132").Value = 5644.875
$ws.Range("M132").Value = -3114.875
$ws.Range("H136").Value = 3350.0454
$ws.Range("J136").Value = 1672
$ws.Range("L136").Value = 5016
$ws.Range("N136").Value = -10116

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4935484.5
$ws.Range("I4").Value = 6169239
$ws.Range("K4").Value = 18507717
$ws.Range("M4").Value = -18507605
$ws.Range("H68").Value = 999
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2997
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2186
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 999
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8991
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4935
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 384936.47
$ws.Range("J86").Value = 398.375
$ws.Range("L86").Value = 1195.125
$ws.Range("N86").Value = -3567.125
$ws.Range("H89").Value = 384936.47
$ws.Range("J89").Value = 398.375
$ws.Range("L89").Value = 3585.375
$ws.Range("N89").Value = -15441.375
$ws.Range("H92").Value = 349.4
$ws.Range("I92").Value = 380
$ws.Range("J92").Value = 329
$ws.Range("K92").Value = 1140
$ws.Range("L92").Value = 987
$ws.Range("M92").Value = 108
$ws.Range("N92").Value = -3483
$ws.Range("H97").Value = 593
$ws.Range("I97").Value = 593
$ws.Range("K97").Value = 1779
$ws.Range("M97").Value = -1283
$ws.Range("H109").Value = 4047.1428
$ws.Range("I109").Value = 3050
$ws.Range("K109").Value = 9150
$ws.Range("M109").Value = -8110
$ws.Range("H122").Value = 1174.2
$ws.Range("J122").Value = 1214.8334
$ws.Range("L122").Value = 10933.5006
$ws.Range("N122").Value = -15833.5006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 206.25
$ws.Range("I2").Value = 206.25
$ws.Range("K2").Value = 206.25
$ws.Range("M2").Value = -93.25
$ws.Range("H41").Value = 6750
$ws.Range("I41").Value = 7500
$ws.Range("K41").Value = 7500
$ws.Range("M41").Value = -7145
$ws.Range("H70").Value = 6763.0586
$ws.Range("I70").Value = 5622.75
$ws.Range("K70").Value = 5622.75
$ws.Range("M70").Value = -5352.75
$ws.Range("H73").Value = 6763.0586
$ws.Range("I73").Value = 5622.75
$ws.Range("K73").Value = 5622.75
$ws.Range("M73").Value = -4686.75
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H132").Value = 4081.4
$ws.Range("I132").Value = 3767.2
$ws.Range("J132").Value = 5024
$ws.Range("K132").Value = 11301.6
$ws.Range("L132").Value = 15072
$ws.Range("M132").Value = -8771.599999999999
$ws.Range("N132").Value = -20132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5225.4614
$ws.Range("J7").Value = 5934.1113
$ws.Range("L7").Value = 5934.1113
$ws.Range("N7").Value = -6158.1113
$ws.Range("H82").Value = 2705.6428
$ws.Range("I82").Value = 2323.3333
$ws.Range("K82").Value = 2323.3333
$ws.Range("M82").Value = -1962.3333
$ws.Range("H85").Value = 2705.6428
$ws.Range("I85").Value = 2323.3333
$ws.Range("K85").Value = 2323.3333
$ws.Range("M85").Value = -1075.3333
$ws.Range("H100").Value = 3494.9285
$ws.Range("I100").Value = 3438
$ws.Range("K100").Value = 3438
$ws.Range("M100").Value = -2897
$ws.Range("H126").Value = 5225.4614
$ws.Range("J126").Value = 5934.1113
$ws.Range("L126").Value = 17802.3339
$ws.Range("N126").Value = -22742.3339
$ws.Range("H132").Value = 3094.1177
$ws.Range("I132").Value = 3190.4
$ws.Range("J132").Value = 2956.5715
$ws.Range("K132").Value = 9571.200000000001
$ws.Range("L132").Value = 8869.7145
$ws.Range("M132").Value = -7041.200000000001
$ws.Range("N132").Value = -13929.7145
$ws.Range("H136").Value = 3522.087
$ws.Range("I136").Value = 2654.5625
$ws.Range("J136").Value = 5505
$ws.Range("K136").Value = 7963.6875
$ws.Range("L136").Value = 16515
$ws.Range("M136").Value = -5413.6875
$ws.Range("N136").Value = -21615

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 43021.2
$ws.Range("I34").Value = 45035.668
$ws.Range("K34").Value = 45035.668
$ws.Range("M34").Value = -44832.668
$ws.Range("H42").Value = 42249.5
$ws.Range("I42").Value = 42249.5
$ws.Range("K42").Value = 42249.5
$ws.Range("M42").Value = -41871.5
$ws.Range("H75").Value = 35000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 35000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 35000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -36872
$ws.Range("H78").Value = 35000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 35000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 105000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -114360
$ws.Range("H88").Value = 39829.668
$ws.Range("J88").Value = 39829.668
$ws.Range("L88").Value = 39829.668
$ws.Range("N88").Value = -40641.668
$ws.Range("H91").Value = 39829.668
$ws.Range("J91").Value = 39829.668
$ws.Range("L91").Value = 39829.668
$ws.Range("N91").Value = -42637.668
$ws.Range("H94").Value = 32500
$ws.Range("J94").Value = 32500
$ws.Range("L94").Value = 32500
$ws.Range("N94").Value = -34302
$ws.Range("H101").Value = 24500
$ws.Range("J101").Value = 24500
$ws.Range("L101").Value = 24500
$ws.Range("N101").Value = -30990
$ws.Range("H107").Value = 1471.3572
$ws.Range("J107").Value = 1838.125
$ws.Range("L107").Value = 5514.375
$ws.Range("N107").Value = -9354.375
$ws.Range("H125").Value = 59999.5
$ws.Range("J125").Value = 59999.5
$ws.Range("L125").Value = 59999.5
$ws.Range("N125").Value = -69839.5
$ws.Range("H132").Value = 2158.36
$ws.Range("I132").Value = 2115
$ws.Range("J132").Value = 2295.6667
$ws.Range("K132").Value = 6345
$ws.Range("L132").Value = 6887.000100000001
$ws.Range("M132").Value = -3815
$ws.Range("N132").Value = -11947.0001
$ws.Range("H136").Value = 5001.25
$ws.Range("I136").Value = 4858.857
$ws.Range("K136").Value = 14576.571
$ws.Range("M136").Value = -12026.571
